$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "hasil_screenshot\screenshot_095849.png"
$ws.Range("C3").Value = "hasil_screenshot\screenshot_095852.png"
$ws.Range("C4").Value = "hasil_screenshot\screenshot_095854.png"
$ws.Range("C5").Value = "hasil_screenshot\screenshot_095856.png"
